$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, matching style of existing header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill I2:J87 with the new data values
$iVals = @(7,8,8,8,8,8,7,7,7,8,7,8,8,8,8,8,8,8,8,9,8,8,7,8,8,7,8,8,8,8,8,8,8,8,9,8,8,8,8,8,8,8,8,8,8,8,8,8,9,8,9,9,8,8,8,8,8,8,8,9,8,8,8,8,8,9,8,8,8,8,7,8,8,8,6,7,8,8,4,6,7,6,5,6,4,3)
$jVals = @(7,8,8,8,8,8,7,7,7,8,8,8,8,8,8,8,8,8,8,9,8,8,8,8,8,7,8,8,8,8,8,8,8,8,9,8,8,8,8,8,8,8,8,8,8,8,8,8,9,8,9,9,8,8,8,8,8,8,8,9,8,8,8,8,8,9,8,8,8,8,7,9,8,8,6,7,8,8,4,6,7,6,5,6,4,3)

for ($r = 0; $r -lt $iVals.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$r]
    $ws.Cells.Item($row, 10).Value = $jVals[$r]
}
